# Auto-generated Excel COM-interop script to apply Leve profit data refresh
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6685.2856
$ws.Range("I18").Value = 7132.8335
$ws.Range("K18").Value = 7132.8335
$ws.Range("M18").Value = -6848.8335

$ws.Range("H32").Value = 2739.8333
$ws.Range("I32").Value = 1280
$ws.Range("J32").Value = 4199.6665
$ws.Range("K32").Value = 1280
$ws.Range("L32").Value = 4199.6665
$ws.Range("M32").Value = -954
$ws.Range("N32").Value = -4851.6665

$ws.Range("H34").Value = 12749.625
$ws.Range("I34").Value = 12856.714
$ws.Range("K34").Value = 12856.714
$ws.Range("M34").Value = -12653.714

$ws.Range("H36").Value = 12749.625
$ws.Range("I36").Value = 12856.714
$ws.Range("K36").Value = 12856.714
$ws.Range("M36").Value = -12141.714

$ws.Range("H100").Value = 2916.3333
$ws.Range("J100").Value = 2916.3333
$ws.Range("L100").Value = 2916.3333
$ws.Range("N100").Value = -3998.3333

$ws.Range("H125").Value = 166677840
$ws.Range("J125").Value = 83349416
$ws.Range("L125").Value = 750144744
$ws.Range("N125").Value = -750149664

$ws.Range("H137").Value = 1649.7
$ws.Range("I137").Value = 1212.125
$ws.Range("K137").Value = 3636.375
$ws.Range("M137").Value = -1086.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3449.875
$ws.Range("I61").Value = 3428.4285
$ws.Range("K61").Value = 3428.4285
$ws.Range("M61").Value = -3216.4285

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H132").Value = 2417.375
$ws.Range("I132").Value = 2532.1428
$ws.Range("K132").Value = 7596.428400000001
$ws.Range("M132").Value = -5066.428400000001

$ws.Range("H136").Value = 3449.875
$ws.Range("I136").Value = 3428.4285
$ws.Range("K136").Value = 10285.2855
$ws.Range("M136").Value = -7735.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2547
$ws.Range("I20").Value = 1308
$ws.Range("K20").Value = 1308
$ws.Range("M20").Value = -1061

$ws.Range("H126").Value = 95999
$ws.Range("J126").Value = 95999
$ws.Range("L126").Value = 95999
$ws.Range("N126").Value = -105879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6795.8
$ws.Range("J16").Value = 7994.5
$ws.Range("L16").Value = 7994.5
$ws.Range("N16").Value = -8568.5

$ws.Range("H62").Value = 5216.3335
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 5759.6
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 5759.6
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -7007.6

$ws.Range("H65").Value = 5216.3335
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 5759.6
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 28798
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -35038

$ws.Range("H113").Value = 6795.8
$ws.Range("J113").Value = 7994.5
$ws.Range("L113").Value = 7994.5
$ws.Range("N113").Value = -12334.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 314
$ws.Range("J62").Value = 314
$ws.Range("L62").Value = 942
$ws.Range("N62").Value = -2314

$ws.Range("H65").Value = 314
$ws.Range("J65").Value = 314
$ws.Range("L65").Value = 2826
$ws.Range("N65").Value = -9690

$ws.Range("H97").Value = 1569
$ws.Range("I97").Value = 2744
$ws.Range("J97").Value = 394
$ws.Range("K97").Value = 8232
$ws.Range("L97").Value = 1182
$ws.Range("M97").Value = -7736
$ws.Range("N97").Value = -2174

$ws.Range("H129").Value = 1252032.1
$ws.Range("J129").Value = 5000991.5
$ws.Range("L129").Value = 15002974.5
$ws.Range("N129").Value = -15012974.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 56400
$ws.Range("J134").Value = 56400
$ws.Range("L134").Value = 169200
$ws.Range("N134").Value = -174270

$ws.Range("H136").Value = 46064.8
$ws.Range("J136").Value = 46064.8
$ws.Range("L136").Value = 138194.4
$ws.Range("N136").Value = -143294.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1499
$ws.Range("I61").Value = 1499.5
$ws.Range("J61").Value = 1498
$ws.Range("K61").Value = 1499.5
$ws.Range("L61").Value = 1498
$ws.Range("M61").Value = -1297.5
$ws.Range("N61").Value = -1902

$ws.Range("H82").Value = 1729.75
$ws.Range("I82").Value = 1566.7333
$ws.Range("J82").Value = 2218.8
$ws.Range("K82").Value = 1566.7333
$ws.Range("L82").Value = 2218.8
$ws.Range("M82").Value = -1205.7333
$ws.Range("N82").Value = -2940.8

$ws.Range("H85").Value = 1729.75
$ws.Range("I85").Value = 1566.7333
$ws.Range("J85").Value = 2218.8
$ws.Range("K85").Value = 1566.7333
$ws.Range("L85").Value = 2218.8
$ws.Range("M85").Value = -318.7333000000001
$ws.Range("N85").Value = -4714.8

$ws.Range("H113").Value = 1499
$ws.Range("I113").Value = 1499.5
$ws.Range("J113").Value = 1498
$ws.Range("K113").Value = 1499.5
$ws.Range("L113").Value = 1498
$ws.Range("M113").Value = 670.5
$ws.Range("N113").Value = -5838

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 74800.75
$ws.Range("I75").Value = 73006
$ws.Range("J75").Value = 75399
$ws.Range("K75").Value = 73006
$ws.Range("L75").Value = 75399
$ws.Range("M75").Value = -72070
$ws.Range("N75").Value = -77271

$ws.Range("H78").Value = 74800.75
$ws.Range("I78").Value = 73006
$ws.Range("J78").Value = 75399
$ws.Range("K78").Value = 219018
$ws.Range("L78").Value = 226197
$ws.Range("M78").Value = -214338
$ws.Range("N78").Value = -235557

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H132").Value = 1999.5
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
